$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "67.838.24"
$ws.Range("E2").Value = "  -1.04%  "

# Row 3
$ws.Range("D3").Value = "3.848.31"
$ws.Range("E3").Value = "  -1.58%  "

# Row 4
$ws.Range("E4").Value = "  +0.01%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "597.61"
$ws.Range("E5").Value = "  -0.88%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "165.82"
$ws.Range("E6").Value = "  +0.76%  "

# Row 7
$ws.Range("D7").Value = "3.845.22"
$ws.Range("E7").Value = "  -1.61%  "

# Row 8
$ws.Range("E8").Value = "  +0.12%  "

# Row 9
$ws.Range("E9").Value = "  -0.20%  "

# Row 10
$ws.Range("E10").Value = "  -0.74%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.33"
$ws.Range("E11").Value = "  -0.39%  "

# Row 12
$ws.Range("E12").Value = "  -0.60%  "

# Row 13
$ws.Range("E13").Value = "  +0.44%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "36.72"
$ws.Range("E14").Value = "  -0.21%  "

# Row 15
$ws.Range("E15").Value = "  -1.44%  "

# Row 16
$ws.Range("D16").Value = "3.857.98"
$ws.Range("E16").Value = "  -0.99%  "

# Row 17
$ws.Range("D17").Value = "67.922.53"

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "18.05"
$ws.Range("E18").Value = "  +6.12%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.33"
$ws.Range("E19").Value = "  -0.78%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.110"
$ws.Range("E20").Value = "  -1.65%  "

# Row 21
$ws.Range("E21").Value = "  -2.69%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "463.26"
$ws.Range("E22").Value = "  -4.29%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.727"
$ws.Range("E23").Value = "  +1.38%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.0000162"
$ws.Range("E24").Value = "  -3.67%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "83.08"
$ws.Range("E25").Value = "  -1.49%  "

# Row 26
$ws.Range("E26").Value = "  +0.19%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.07"
$ws.Range("E27").Value = "  +0.93%  "

# Row 28
$ws.Range("E28").Value = "  -0.02%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.95"
$ws.Range("E29").Value = "  -1.14%  "

# Row 30
$ws.Range("E30").Value = "  +0.38%  "

# Row 31
$ws.Range("D31").Value = "3.999.67"
$ws.Range("E31").Value = "  -1.45%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.71"
$ws.Range("E32").Value = "  -1.32%  "

# Row 33
$ws.Range("E33").Value = "  -2.36%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "30.96"
$ws.Range("E34").Value = "  -3.11%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "9.28"
$ws.Range("E35").Value = "  +1.77%  "

# Row 36
$ws.Range("D36").Value = "3.826.24"
$ws.Range("E36").Value = "  -0.69%  "

# Row 37
$ws.Range("E37").Value = "  -2.24%  "

# Row 38
$ws.Range("E38").Value = "  -1.58%  "

# Row 39
$ws.Range("E39").Value = "  -0.19%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.87"
$ws.Range("E40").Value = "  +0.05%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.25"
$ws.Range("E41").Value = "  +5.74%  "

# Row 42
$ws.Range("E42").Value = "  +0.07%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.311"
$ws.Range("E43").Value = "  -1.52%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "426.09"
$ws.Range("E44").Value = "  -1.61%  "

# Row 45
$ws.Range("E45").Value = "  -0.34%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "47.12"
$ws.Range("E47").Value = "  -2.80%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "8.47"
$ws.Range("E48").Value = "  +0.78%  "

# Row 49
$ws.Range("B49").Value = "FLOKI"
$ws.Range("C49").Value = "https://coinranking.com/coin/fmHk13Rqw+floki-floki"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.000273"
$ws.Range("E49").Value = "  +3.69%  "

# Row 50
$ws.Range("B50").Value = "Monero"
$ws.Range("C50").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "143.62"
$ws.Range("E50").Value = "  +1.23%  "

# Row 51
$ws.Range("B51").Value = "Arweave"
$ws.Range("C51").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "40.30"
$ws.Range("E51").Value = "  +3.12%  "
